## Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
# Reformats the "Source" block of the Mexico Summary sheet:
#   - splits the old "INEGI sobre el tema..." paragraph onto its own row
#   - turns the hyperlinked INEGI URL into a plain (non-hyperlinked) text row
#   - moves the "INEGI" citation label down one row and repeats it as a plain source line

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The long descriptive blurb that used to live in A38
$inegiBlurb = "INEGI sobre el tema:`nResumen de los resultados de los Censos Económicos 2009; Micro, pequeña, mediana, y  gran empresa. Estratificación de los establecimientos, 2009"

# The INEGI source URL that used to live in A39 (as a hyperlink)
$inegiUrl = "http://www.inegi.org.mx/prod_serv/contenidos/espanol/bvinegi/productos/censos/economicos/2009/comercio/micro_peque_media/Mono_Micro_peque_mediana.pdf"

# Drop the hyperlink on A39 entirely (it becomes plain text further down)
$ws.Hyperlinks.Delete()

# Row 38 becomes blank (text moves down to row 39)
$ws.Range("A38").Value = ""

# Row 39: the long INEGI paragraph, styled like the other italic "source" rows
$ws.Range("A39").Value = $inegiBlurb
$ws.Range("A39").Font().Italic = $true

# Row 40 stays blank (unchanged)

# Row 41 (new): the INEGI URL as plain text, no hyperlink, same italic "source" style
$ws.Range("A41").Value = $inegiUrl
$ws.Range("A41").Font().Italic = $true

# Re-apply the same italic styling to the rows that already looked this way so
# every "source" row ends up sharing one consistent style
$ws.Range("A37").Font().Italic = $true
$ws.Range("A40").Font().Italic = $true

# Row 43 (old "INEGI" title line) is removed entirely - it now lives on row 44
$ws.Range("A43").Clear()

# Row 44: "INEGI" title line (bold "title" style), moved down from row 43
$ws.Range("A44").Value = "INEGI"
$ws.Range("A44").Font().Bold = $true

# Row 45 (new): a plain "INEGI" source line (replaces the old long legal citation)
$ws.Range("A45").Value = "INEGI"
$ws.Range("A45").Font().Italic = $true

Write-Host "Edit applied"
